# Add data for 2022-11-17
# Updates the "through November 08" running total column to "through November 09",
# incrementing counts for neighborhoods that had an additional carjacking recorded
# on 2022-11-09 (and adding brand-new entries where a neighborhood/month cell
# previously had no incidents).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet tab and update the "November 2022" column header text to
# reflect the new as-of date.
$ws.Name = "Through 2022-11-09"
$ws.Range("B1").Value = "November 2022 (through November 09)"

# --- Cells whose existing value increments by 1 ---
$valueUpdates = @{
    "M2"   = 3
    "CA3"  = 2
    "B5"   = 4
    "AI5"  = 2
    "X14"  = 6
    "X17"  = 2
    "M21"  = 6
    "M24"  = 3
    "B25"  = 2
    "M26"  = 3
    "M27"  = 2
    "BE33" = 2
    "BE41" = 2
    "M50"  = 2
    "X74"  = 2
    "X76"  = 3
    "M84"  = 2
}

foreach ($addr in $valueUpdates.Keys) {
    $ws.Range($addr).Value = $valueUpdates[$addr]
}

# --- Cells that are brand new (previously empty, now hold a count of 1) ---
$newCells = @(
    "M12",
    "AT13",
    "BP13",
    "AT14",
    "BP14",
    "CA14",
    "M17",
    "B18",
    "X35",
    "AT41",
    "X43",
    "M57",
    "X64",
    "BP64",
    "BE72"
)

foreach ($addr in $newCells) {
    $ws.Range($addr).Value = 1
}
